# Apply the cyclic re-shuffle of data for rows 5-8 as described in the diff.
# The underlying species records (columns A,B,D,E,F,G,H,Q,R) were rotated
# between the existing rows while the remaining, identical metadata columns
# (C,I,P,S,T,U,V,W,Y,AA,AD,AE,AG,AT,AW,AX,AY) stay exactly where they were.
#
# Net effect (old row -> new row):
#   old row 6 -> new row 5
#   old row 8 -> new row 6
#   old row 5 -> new row 7
#   old row 7 -> new row 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose values move with the record.
$cols = @("A","B","D","E","F","G","H","Q","R")

# Capture the current ("before") values for rows 5-8 for each of those columns.
# Use Value2 (rather than Value) to avoid COM Variant wrapping issues.
$before = @{}
foreach ($r in 5..8) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Mapping of new row -> source (old) row.
$mapping = @{ 5 = 6; 6 = 8; 7 = 5; 8 = 7 }

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $before[$oldRow][$c]
    }
}
